$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "Invincibilità momentan" + bookmark(_GoBack) + "ea" -> single run
# "Invincibilità momentanea" with the bookmark removed. Find/Replace
# across the whole visible phrase merges the split runs into one and
# drops the bookmark that sat in the middle of it.
$d.Content.Find.Execute("Invincibilità momentanea", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Invincibilità momentanea", 2)

# --- Change 2 ---------------------------------------------------------
# The first of the two empty paragraphs right after the "Sputano" fuoco
# bullet gets new text, plus a (collapsed) _GoBack bookmark right after
# the text (before the paragraph mark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Sputano*fuoco*particellari*") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Fill in the text. A trailing placeholder character is appended so the
# collapsed bookmark we add next sits one character before the (buggy)
# "end of paragraph text" position; the placeholder is then deleted,
# leaving the bookmark correctly collapsed right after the real text.
$target.Range.Text = "Come checkpoint si potrebbe mettere un cancelloZ"

$filled = $target.Range
$bmPos = $filled.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($filled.End - 2, $filled.End - 1)
$placeholder.Delete()
